# This script re-orders the per-row "observation" values (Fecha, Calidad, Volumen,
# Precio minimo/maximo/promedio, Unidad de comercializacion, Origen, Precio $/Kg,
# Kg o Unidades -- columns D and I..Q) across rows 2..104, while leaving the
# "constant" columns (A,B,C,E,F,G,H,R) untouched, since each of those columns
# holds the same repeated value for every data row in this sheet.
#
# Concretely: for output row r (2..104), the new D/I/J/K/L/M/N/O/P/Q values are
# copied from the *original* (pre-edit) row given by $sourceRow[r-2].

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 104

# sourceRow[i] = original row number whose D..Q payload should end up in row (firstRow + i)
$sourceRow = @(87,38,71,84,31,35,80,26,91,5,24,63,65,22,88,30,93,13,14,34,99,104,103,10,51,81,54,62,78,90,75,76,82,11,89,27,86,101,9,4,43,58,59,70,94,40,52,6,2,44,45,36,56,55,32,77,50,25,37,3,49,18,66,8,102,41,46,95,72,60,61,92,21,12,17,39,85,64,67,53,29,68,47,57,74,42,97,15,79,48,73,98,16,83,19,20,28,23,69,96,7,100,33)

$cols = @("D","I","J","K","L","M","N","O","P","Q")

# 1) Snapshot every value we might need to read, BEFORE any writes happen,
#    so that writes to one row never clobber data still needed as a source
#    for another row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write the permuted values back out.
for ($i = 0; $i -lt $sourceRow.Length; $i++) {
    $destRow = $firstRow + $i
    $srcRow = $sourceRow[$i]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $srcVals[$c]
    }
}
